# Generate Report for Handoff
# Swaps the order of the "1c57ed33..." and "95228697..." file rows on every
# sheet (the 95228697 entry now appears first / row 2, the 1c57ed33 entry
# now appears second / row 3) and marks the 1c57ed33 entry as freshly
# "Ready for handoff" with updated handoff timestamps.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet 1: Overview
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("Overview")

$ws1.Range("A2").Value = "95228697-62d2-43dd-8899-dfd6bcdb101a.md"
$ws1.Range("B2").Value = "Handed back: in sync with en-US"
$ws1.Range("C2").Value = "Handed back: in sync with en-US"
$ws1.Range("D2").Value = "2016-03-23 02:46:00"

$ws1.Range("A3").Value = "1c57ed33-8833-4edc-beb6-3d3e3b3326da.md"
$ws1.Range("B3").Value = "Ready for handoff"
$ws1.Range("C3").Value = "Ready for handoff"
$ws1.Range("D3").Value = "2016-03-23 02:47:25"

# Recreate hyperlinks, keeping the original Address order (so relationship
# ids / targets stay exactly as before) but with the new display text.
$ws1.Hyperlinks.Delete()
$ws1.Hyperlinks.Add($ws1.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/e318d4c7cd246b8c9faa9243e0dabf29e7397050/e2e/1c57ed33-8833-4edc-beb6-3d3e3b3326da.md", "", "", "95228697-62d2-43dd-8899-dfd6bcdb101a.md") | Out-Null
$ws1.Hyperlinks.Add($ws1.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/e318d4c7cd246b8c9faa9243e0dabf29e7397050/e2e/95228697-62d2-43dd-8899-dfd6bcdb101a.md", "", "", "1c57ed33-8833-4edc-beb6-3d3e3b3326da.md") | Out-Null

# ---------------------------------------------------------------------------
# Sheet 2: zh-cn
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("zh-cn")

$ws2.Range("A2").Value = "95228697-62d2-43dd-8899-dfd6bcdb101a.md"
$ws2.Range("B2").Value = ".md"
$ws2.Range("C2").Value = "Handed back: in sync with en-US"
$ws2.Range("D2").Value = "95228697-62d2-43dd-8899-dfd6bcdb101a.5748eb9dee8d2a25dbff5c838c5e4c76ddae178e.zh-cn.xlf"
$ws2.Range("E2").Value = "2016-03-23 02:45:55"
$ws2.Range("F2").Value = "95228697-62d2-43dd-8899-dfd6bcdb101a.md"
$ws2.Range("G2").Value = "95228697-62d2-43dd-8899-dfd6bcdb101a.5748eb9dee8d2a25dbff5c838c5e4c76ddae178e.zh-cn.xlf"
$ws2.Range("H2").Value = "2016-03-23 02:46:30"
$ws2.Range("J2").Value = "Include"

$ws2.Range("A3").Value = "1c57ed33-8833-4edc-beb6-3d3e3b3326da.md"
$ws2.Range("B3").Value = ".md"
$ws2.Range("C3").Value = "Ready for handoff"
$ws2.Range("D3").Value = "1c57ed33-8833-4edc-beb6-3d3e3b3326da.4192ab7586250a185305b5734ef0f6ba9ae63bc1.zh-cn.xlf"
$ws2.Range("E3").Value = "2016-03-23 02:47:20"
$ws2.Range("F3").Value = "1c57ed33-8833-4edc-beb6-3d3e3b3326da.md"
$ws2.Range("G3").Value = "1c57ed33-8833-4edc-beb6-3d3e3b3326da.4192ab7586250a185305b5734ef0f6ba9ae63bc1.zh-cn.xlf"
$ws2.Range("H3").Value = "2016-03-23 02:46:30"
$ws2.Range("J3").Value = "Include"

$ws2.Hyperlinks.Delete()
$ws2.Hyperlinks.Add($ws2.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/e318d4c7cd246b8c9faa9243e0dabf29e7397050/e2e/1c57ed33-8833-4edc-beb6-3d3e3b3326da.md", "", "", "95228697-62d2-43dd-8899-dfd6bcdb101a.md") | Out-Null
$ws2.Hyperlinks.Add($ws2.Range("D2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/4937c1ed291fd931e91e56237edd6e48561f2733/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/1c57ed33-8833-4edc-beb6-3d3e3b3326da.4192ab7586250a185305b5734ef0f6ba9ae63bc1.zh-cn.xlf", "", "", "95228697-62d2-43dd-8899-dfd6bcdb101a.5748eb9dee8d2a25dbff5c838c5e4c76ddae178e.zh-cn.xlf") | Out-Null
$ws2.Hyperlinks.Add($ws2.Range("F2"), "https://github.com/OpenLocalizationTestOrg/oltest.zh-cn/blob/4796d1d716c6f271e8f0322b5115bf6acb207947/e2e/1c57ed33-8833-4edc-beb6-3d3e3b3326da.md", "", "", "95228697-62d2-43dd-8899-dfd6bcdb101a.md") | Out-Null
$ws2.Hyperlinks.Add($ws2.Range("G2"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/10ebd387ca3f93aca1c3bf5bf3058a13500be9c1/ol-handback/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/1c57ed33-8833-4edc-beb6-3d3e3b3326da.4192ab7586250a185305b5734ef0f6ba9ae63bc1.zh-cn.xlf", "", "", "95228697-62d2-43dd-8899-dfd6bcdb101a.5748eb9dee8d2a25dbff5c838c5e4c76ddae178e.zh-cn.xlf") | Out-Null
$ws2.Hyperlinks.Add($ws2.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/e318d4c7cd246b8c9faa9243e0dabf29e7397050/e2e/95228697-62d2-43dd-8899-dfd6bcdb101a.md", "", "", "1c57ed33-8833-4edc-beb6-3d3e3b3326da.md") | Out-Null
$ws2.Hyperlinks.Add($ws2.Range("D3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/4937c1ed291fd931e91e56237edd6e48561f2733/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/95228697-62d2-43dd-8899-dfd6bcdb101a.5748eb9dee8d2a25dbff5c838c5e4c76ddae178e.zh-cn.xlf", "", "", "1c57ed33-8833-4edc-beb6-3d3e3b3326da.4192ab7586250a185305b5734ef0f6ba9ae63bc1.zh-cn.xlf") | Out-Null
$ws2.Hyperlinks.Add($ws2.Range("F3"), "https://github.com/OpenLocalizationTestOrg/oltest.zh-cn/blob/4796d1d716c6f271e8f0322b5115bf6acb207947/e2e/95228697-62d2-43dd-8899-dfd6bcdb101a.md", "", "", "1c57ed33-8833-4edc-beb6-3d3e3b3326da.md") | Out-Null
$ws2.Hyperlinks.Add($ws2.Range("G3"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/10ebd387ca3f93aca1c3bf5bf3058a13500be9c1/ol-handback/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/95228697-62d2-43dd-8899-dfd6bcdb101a.5748eb9dee8d2a25dbff5c838c5e4c76ddae178e.zh-cn.xlf", "", "", "1c57ed33-8833-4edc-beb6-3d3e3b3326da.4192ab7586250a185305b5734ef0f6ba9ae63bc1.zh-cn.xlf") | Out-Null

# ---------------------------------------------------------------------------
# Sheet 3: de-de
# ---------------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("de-de")

$ws3.Range("A2").Value = "95228697-62d2-43dd-8899-dfd6bcdb101a.md"
$ws3.Range("B2").Value = ".md"
$ws3.Range("C2").Value = "Handed back: in sync with en-US"
$ws3.Range("D2").Value = "95228697-62d2-43dd-8899-dfd6bcdb101a.5748eb9dee8d2a25dbff5c838c5e4c76ddae178e.de-de.xlf"
$ws3.Range("E2").Value = "2016-03-23 02:46:00"
$ws3.Range("F2").Value = "95228697-62d2-43dd-8899-dfd6bcdb101a.md"
$ws3.Range("G2").Value = "95228697-62d2-43dd-8899-dfd6bcdb101a.5748eb9dee8d2a25dbff5c838c5e4c76ddae178e.de-de.xlf"
$ws3.Range("H2").Value = "2016-03-23 02:46:38"
$ws3.Range("J2").Value = "Include"

$ws3.Range("A3").Value = "1c57ed33-8833-4edc-beb6-3d3e3b3326da.md"
$ws3.Range("B3").Value = ".md"
$ws3.Range("C3").Value = "Ready for handoff"
$ws3.Range("D3").Value = "1c57ed33-8833-4edc-beb6-3d3e3b3326da.4192ab7586250a185305b5734ef0f6ba9ae63bc1.de-de.xlf"
$ws3.Range("E3").Value = "2016-03-23 02:47:25"
$ws3.Range("F3").Value = "1c57ed33-8833-4edc-beb6-3d3e3b3326da.md"
$ws3.Range("G3").Value = "1c57ed33-8833-4edc-beb6-3d3e3b3326da.4192ab7586250a185305b5734ef0f6ba9ae63bc1.de-de.xlf"
$ws3.Range("H3").Value = "2016-03-23 02:46:38"
$ws3.Range("J3").Value = "Include"

$ws3.Hyperlinks.Delete()
$ws3.Hyperlinks.Add($ws3.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/e318d4c7cd246b8c9faa9243e0dabf29e7397050/e2e/1c57ed33-8833-4edc-beb6-3d3e3b3326da.md", "", "", "95228697-62d2-43dd-8899-dfd6bcdb101a.md") | Out-Null
$ws3.Hyperlinks.Add($ws3.Range("D2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/03a114baa5022207396ff2f920ac3b8278ce6d3b/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/1c57ed33-8833-4edc-beb6-3d3e3b3326da.4192ab7586250a185305b5734ef0f6ba9ae63bc1.de-de.xlf", "", "", "95228697-62d2-43dd-8899-dfd6bcdb101a.5748eb9dee8d2a25dbff5c838c5e4c76ddae178e.de-de.xlf") | Out-Null
$ws3.Hyperlinks.Add($ws3.Range("F2"), "https://github.com/OpenLocalizationTestOrg/oltest.de-de/blob/d5ce135141729085bb881b294650bb24cbfc545a/e2e/1c57ed33-8833-4edc-beb6-3d3e3b3326da.md", "", "", "95228697-62d2-43dd-8899-dfd6bcdb101a.md") | Out-Null
$ws3.Hyperlinks.Add($ws3.Range("G2"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/815c2ee04719b328c02744438f9d01b92eb22288/ol-handback/OpenLocalizationTestOrg/oltest.de-de/ci/ht/1c57ed33-8833-4edc-beb6-3d3e3b3326da.4192ab7586250a185305b5734ef0f6ba9ae63bc1.de-de.xlf", "", "", "95228697-62d2-43dd-8899-dfd6bcdb101a.5748eb9dee8d2a25dbff5c838c5e4c76ddae178e.de-de.xlf") | Out-Null
$ws3.Hyperlinks.Add($ws3.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/e318d4c7cd246b8c9faa9243e0dabf29e7397050/e2e/95228697-62d2-43dd-8899-dfd6bcdb101a.md", "", "", "1c57ed33-8833-4edc-beb6-3d3e3b3326da.md") | Out-Null
$ws3.Hyperlinks.Add($ws3.Range("D3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/03a114baa5022207396ff2f920ac3b8278ce6d3b/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/95228697-62d2-43dd-8899-dfd6bcdb101a.5748eb9dee8d2a25dbff5c838c5e4c76ddae178e.de-de.xlf", "", "", "1c57ed33-8833-4edc-beb6-3d3e3b3326da.4192ab7586250a185305b5734ef0f6ba9ae63bc1.de-de.xlf") | Out-Null
$ws3.Hyperlinks.Add($ws3.Range("F3"), "https://github.com/OpenLocalizationTestOrg/oltest.de-de/blob/d5ce135141729085bb881b294650bb24cbfc545a/e2e/95228697-62d2-43dd-8899-dfd6bcdb101a.md", "", "", "1c57ed33-8833-4edc-beb6-3d3e3b3326da.md") | Out-Null
$ws3.Hyperlinks.Add($ws3.Range("G3"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/815c2ee04719b328c02744438f9d01b92eb22288/ol-handback/OpenLocalizationTestOrg/oltest.de-de/ci/ht/95228697-62d2-43dd-8899-dfd6bcdb101a.5748eb9dee8d2a25dbff5c838c5e4c76ddae178e.de-de.xlf", "", "", "1c57ed33-8833-4edc-beb6-3d3e3b3326da.4192ab7586250a185305b5734ef0f6ba9ae63bc1.de-de.xlf") | Out-Null

$wb.Save()
